$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.526893019676208
$ws.Range("B1").Value = 2.086515188217163
$ws.Range("C1").Value = 3.227139711380005
$ws.Range("D1").Value = 4.830220222473145
$ws.Range("E1").Value = 0.8482956290245056
